# Apply "new chain_detailsresults scneario 1 and 7" edit.
#
# Summary of changes:
#  - Chain_Details: add a "desired_gradient" column (G) and two new rows
#    (scenario_id 7) for Chain 1 (popBased, gradient 0.2) and Chain 2
#    (popBased / uniform).
#  - Outbreaks: add a matching row (scenario_id 7) with the outbreak
#    parameters for that new scenario.
#  - Explanation of Options: remove the now-unused placeholder cells A8/B8.
#  - Selection / active-sheet bookkeeping to match the saved UI state.

$wb = $excel.ActiveWorkbook

$wsChain = $wb.Worksheets.Item("Chain_Details")
$wsPop   = $wb.Worksheets.Item("Population")
$wsOut   = $wb.Worksheets.Item("Outbreaks")
$wsExpl  = $wb.Worksheets.Item("Explanation of Options")

# ---------------------------------------------------------------------
# 1. Chain_Details (sheet1) - new column G "desired_gradient" + rows 4/5
# ---------------------------------------------------------------------

# Header cell G1, matching the style of the other header cells (A1:F1).
$wsChain.Range("A1").Copy()
$wsChain.Range("G1").PasteSpecial(-4122)
$wsChain.Range("G1").Value = "desired_gradient"

# G2 / G3 are "not applicable" for the existing rows - give them the same
# grey placeholder formatting used elsewhere in the workbook (e.g.
# Population!D2) for cells where a parameter doesn't apply.
$wsPop.Range("D2").Copy()
$wsChain.Range("G2").PasteSpecial(-4122)
$wsPop.Range("D2").Copy()
$wsChain.Range("G3").PasteSpecial(-4122)

# New row 4: scenario 7 / Chain 1 / popBased / popBased, with a gradient.
$wsChain.Range("A2:F2").Copy()
$wsChain.Range("A4").PasteSpecial(-4122)
$wsChain.Range("F2").Copy()
$wsChain.Range("G4").PasteSpecial(-4122)
$wsChain.Range("A4").Value = 7
$wsChain.Range("B4").Value = "Chain 1"
$wsChain.Range("C4").Value = 3
$wsChain.Range("D4").Value = "popBased"
$wsChain.Range("E4").Value = "popBased"
$wsChain.Range("F4").Value = 1000
$wsChain.Range("G4").Value = 0.2

# New row 5: scenario 7 / Chain 2 / popBased / uniform, no gradient.
$wsChain.Range("A3:F3").Copy()
$wsChain.Range("A5").PasteSpecial(-4122)
$wsPop.Range("D2").Copy()
$wsChain.Range("G5").PasteSpecial(-4122)
$wsChain.Range("A5").Value = 7
$wsChain.Range("B5").Value = "Chain 2"
$wsChain.Range("C5").Value = 3
$wsChain.Range("D5").Value = "popBased"
$wsChain.Range("E5").Value = "uniform"
$wsChain.Range("F5").Value = 1000

# Widen the new column to fit its header text.
$wsChain.Columns.Item(7).ColumnWidth = 13.75

# ---------------------------------------------------------------------
# 2. Outbreaks (sheet3) - add matching scenario-7 outbreak parameters
# ---------------------------------------------------------------------

$wsOut.Range("A2:E2").Copy()
$wsOut.Range("A3").PasteSpecial(-4122)
$wsOut.Range("A3").Value = 7
$wsOut.Range("B3").Value = 0.2
$wsOut.Range("C3").Value = 0.001
$wsOut.Range("D3").Value = 5
$wsOut.Range("E3").Value = 1

# ---------------------------------------------------------------------
# 3. Explanation of Options (sheet4) - drop the stray A8/B8 cells
# ---------------------------------------------------------------------

$wsExpl.Range("A8").Clear()
$wsExpl.Range("B8").Clear()

# ---------------------------------------------------------------------
# 4. Selection / active sheet bookkeeping
# ---------------------------------------------------------------------

$wsPop.Range("D1:D17").Select()
$wsOut.Range("B12").Select()
$wsChain.Range("G14").Select()
$wsChain.Activate()
